$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 95
$ws1.Range("F4").Value = 607
$ws1.Range("F5").Value = 139
$ws1.Range("F6").Value = 9184
$ws1.Range("F9").Value = 1174
$ws1.Range("F10").Value = 1060
$ws1.Range("F11").Value = 139
$ws1.Range("F14").Value = 247
$ws1.Range("F15").Value = 360
$ws1.Range("F18").Value = 1194

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 95
$ws4.Range("F6").Value = 607
$ws4.Range("F7").Value = 139
$ws4.Range("F8").Value = 9184
$ws4.Range("F11").Value = 1174
$ws4.Range("F12").Value = 1060
$ws4.Range("F13").Value = 139
$ws4.Range("F16").Value = 247
$ws4.Range("F17").Value = 360
$ws4.Range("F20").Value = 1194
